# Re-curate the "grandes-grupos", "sexo" and "provincia" columns so that
# they are treated as iaest measures instead of dimensions, and drop the
# now-obsolete mapping-file row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E: grandes-grupos
$ws.Range("E2").Value = "iaest-measure:grandes-grupos"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"

# Column I: sexo
$ws.Range("I2").Value = "iaest-measure:sexo"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"

# Column K: provincia
$ws.Range("K2").Value = "iaest-measure:provincia"
$ws.Range("K3").Value = "medida"
$ws.Range("K4").Value = "xsd:int"

# Row 5 held the mapping-file references for the old dimensions; it is no
# longer needed now that grandes-grupos/sexo/provincia are measures.
$ws.Rows.Item(5).Delete()
